# "Add data for 2024-11-20": append a new column BU to the single
# worksheet ("合成確率"), one cell per existing row (1-53). Row 1 holds
# the new date label; rows 2-53 hold that day's numeric value, styled
# with whichever of the sheet's three cell styles (plain / yellow /
# light-blue) the row's other day-columns already use.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column BU (73) gets the same width as every other date column.
$ws.Columns.Item(73).ColumnWidth = 11.17

$styleAnchor1 = $ws.Range("A2")   # s="1": メイリオ, no fill
$styleAnchor2 = $ws.Range("H8")   # s="2": メイリオ, yellow fill
$styleAnchor3 = $ws.Range("B7")   # s="3": メイリオ, light-blue fill

$styleAnchor1.Copy()
$ws.Range("BU1").PasteSpecial(-4122)
$ws.Range("BU1").NumberFormat = "@"   # keep literal text, not an auto-converted date
$ws.Range("BU1").Value = "2024/11/20"
$ws.Range("BU1").NumberFormat = "General"
$styleAnchor1.Copy()                        # re-apply style (NumberFormat churn above can touch it)
$ws.Range("BU1").PasteSpecial(-4122)

$styleAnchor1.Copy()
$ws.Range("BU2").PasteSpecial(-4122)
$ws.Range("BU2").Value = 183.1

$styleAnchor1.Copy()
$ws.Range("BU3").PasteSpecial(-4122)
$ws.Range("BU3").Value = 176.3

$styleAnchor1.Copy()
$ws.Range("BU4").PasteSpecial(-4122)
$ws.Range("BU4").Value = 145.6

$styleAnchor1.Copy()
$ws.Range("BU5").PasteSpecial(-4122)
$ws.Range("BU5").Value = 202.6

$styleAnchor1.Copy()
$ws.Range("BU6").PasteSpecial(-4122)
$ws.Range("BU6").Value = 202.8

$styleAnchor3.Copy()
$ws.Range("BU7").PasteSpecial(-4122)
$ws.Range("BU7").Value = 130.3

$styleAnchor2.Copy()
$ws.Range("BU8").PasteSpecial(-4122)
$ws.Range("BU8").Value = 123.2

$styleAnchor1.Copy()
$ws.Range("BU9").PasteSpecial(-4122)
$ws.Range("BU9").Value = 145.8

$styleAnchor1.Copy()
$ws.Range("BU10").PasteSpecial(-4122)
$ws.Range("BU10").Value = 160.1

$styleAnchor1.Copy()
$ws.Range("BU11").PasteSpecial(-4122)
$ws.Range("BU11").Value = 155.2

$styleAnchor1.Copy()
$ws.Range("BU12").PasteSpecial(-4122)
$ws.Range("BU12").Value = 171.3

$styleAnchor3.Copy()
$ws.Range("BU13").PasteSpecial(-4122)
$ws.Range("BU13").Value = 136.4

$styleAnchor3.Copy()
$ws.Range("BU14").PasteSpecial(-4122)
$ws.Range("BU14").Value = 127.6

$styleAnchor1.Copy()
$ws.Range("BU15").PasteSpecial(-4122)
$ws.Range("BU15").Value = 142.3

$styleAnchor1.Copy()
$ws.Range("BU16").PasteSpecial(-4122)
$ws.Range("BU16").Value = 142.5

$styleAnchor1.Copy()
$ws.Range("BU17").PasteSpecial(-4122)
$ws.Range("BU17").Value = 174.2

$styleAnchor1.Copy()
$ws.Range("BU18").PasteSpecial(-4122)
$ws.Range("BU18").Value = 168.1

$styleAnchor1.Copy()
$ws.Range("BU19").PasteSpecial(-4122)
$ws.Range("BU19").Value = 185.1

$styleAnchor1.Copy()
$ws.Range("BU20").PasteSpecial(-4122)
$ws.Range("BU20").Value = 255.7

$styleAnchor1.Copy()
$ws.Range("BU21").PasteSpecial(-4122)
$ws.Range("BU21").Value = 161.5

$styleAnchor1.Copy()
$ws.Range("BU22").PasteSpecial(-4122)
$ws.Range("BU22").Value = 140.2

$styleAnchor1.Copy()
$ws.Range("BU23").PasteSpecial(-4122)
$ws.Range("BU23").Value = 153.6

$styleAnchor1.Copy()
$ws.Range("BU24").PasteSpecial(-4122)
$ws.Range("BU24").Value = 145

$styleAnchor1.Copy()
$ws.Range("BU25").PasteSpecial(-4122)
$ws.Range("BU25").Value = 191.5

$styleAnchor1.Copy()
$ws.Range("BU26").PasteSpecial(-4122)
$ws.Range("BU26").Value = 176.3

$styleAnchor1.Copy()
$ws.Range("BU27").PasteSpecial(-4122)
$ws.Range("BU27").Value = 174.2

$styleAnchor1.Copy()
$ws.Range("BU28").PasteSpecial(-4122)
$ws.Range("BU28").Value = 175.6

$styleAnchor1.Copy()
$ws.Range("BU29").PasteSpecial(-4122)
$ws.Range("BU29").Value = 167.6

$styleAnchor1.Copy()
$ws.Range("BU30").PasteSpecial(-4122)
$ws.Range("BU30").Value = 190.1

$styleAnchor1.Copy()
$ws.Range("BU31").PasteSpecial(-4122)
$ws.Range("BU31").Value = 187.7

$styleAnchor2.Copy()
$ws.Range("BU32").PasteSpecial(-4122)
$ws.Range("BU32").Value = 123.9

$styleAnchor1.Copy()
$ws.Range("BU33").PasteSpecial(-4122)
$ws.Range("BU33").Value = 184.2

$styleAnchor1.Copy()
$ws.Range("BU34").PasteSpecial(-4122)
$ws.Range("BU34").Value = 216

$styleAnchor1.Copy()
$ws.Range("BU35").PasteSpecial(-4122)
$ws.Range("BU35").Value = 151.2

$styleAnchor1.Copy()
$ws.Range("BU36").PasteSpecial(-4122)
$ws.Range("BU36").Value = 155

$styleAnchor1.Copy()
$ws.Range("BU37").PasteSpecial(-4122)
$ws.Range("BU37").Value = 155.2

$styleAnchor3.Copy()
$ws.Range("BU38").PasteSpecial(-4122)
$ws.Range("BU38").Value = 136.7

$styleAnchor1.Copy()
$ws.Range("BU39").PasteSpecial(-4122)
$ws.Range("BU39").Value = 155.1

$styleAnchor3.Copy()
$ws.Range("BU40").PasteSpecial(-4122)
$ws.Range("BU40").Value = 132.2

$styleAnchor1.Copy()
$ws.Range("BU41").PasteSpecial(-4122)
$ws.Range("BU41").Value = 156.3

$styleAnchor3.Copy()
$ws.Range("BU42").PasteSpecial(-4122)
$ws.Range("BU42").Value = 134.7

$styleAnchor3.Copy()
$ws.Range("BU43").PasteSpecial(-4122)
$ws.Range("BU43").Value = 134.6

$styleAnchor1.Copy()
$ws.Range("BU44").PasteSpecial(-4122)
$ws.Range("BU44").Value = 147.9

$styleAnchor1.Copy()
$ws.Range("BU45").PasteSpecial(-4122)
$ws.Range("BU45").Value = 206

$styleAnchor1.Copy()
$ws.Range("BU46").PasteSpecial(-4122)
$ws.Range("BU46").Value = 170

$styleAnchor1.Copy()
$ws.Range("BU47").PasteSpecial(-4122)
$ws.Range("BU47").Value = 172.1

$styleAnchor1.Copy()
$ws.Range("BU48").PasteSpecial(-4122)
$ws.Range("BU48").Value = 157.2

$styleAnchor3.Copy()
$ws.Range("BU49").PasteSpecial(-4122)
$ws.Range("BU49").Value = 132.9

$styleAnchor3.Copy()
$ws.Range("BU50").PasteSpecial(-4122)
$ws.Range("BU50").Value = 137

$styleAnchor1.Copy()
$ws.Range("BU51").PasteSpecial(-4122)
$ws.Range("BU51").Value = 190

$styleAnchor1.Copy()
$ws.Range("BU52").PasteSpecial(-4122)
$ws.Range("BU52").Value = 169.5

$styleAnchor1.Copy()
$ws.Range("BU53").PasteSpecial(-4122)
$ws.Range("BU53").Value = 173.1

$excel.CutCopyMode = $false